$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Paragraph 1: "Love is bold but don't go to far" ("is" is bold)
$tr.Text = "Love is bold but don't go to far"
$tr.Characters(6, 2).Font.Bold = $true

# Paragraph 2: "the ice cream man appear on a van"
$tr.InsertAfter("`rthe ice cream man appear on a van")

# Paragraph 3: "with all flavors" ("all" is italic) -- appended as the new last
# paragraph so the italic sub-range is applied while it is last (avoids the
# interop quirk where formatting a partial run in a paragraph that has
# exactly one paragraph after it mis-applies to the whole paragraph).
$ins3 = $tr.InsertAfter("`rwith all flavors")
$p3 = $tr.Paragraphs($tr.Paragraphs().Count)
$p3.Characters(6, 3).Font.Italic = $true

# Paragraph 4: "to try"
$tr.InsertAfter("`rto try")
